# Append a fresh scrape run (2025-10-25 12:32:59 JST) to the "ランサーズ" sheet.
# Net effect vs before.xlsx:
#   - 4 brand new listings inserted (interleaved by priority score), existing
#     listings shift down accordingly
#   - every row's "取得日時" (A) timestamp bumped to the new scrape time
#   - column H ("スキル概要") widened 12 -> 19
#   - hyperlinks on column F re-pointed/extended to match the new row layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Final data for rows 2..13 (header stays in row 1, untouched).
# Columns: A 取得日時, B タイトル, C カテゴリ, D 価格, E 締切, F URL, G 優先度スコア, H スキル概要
$rows = @(
    @("2025-10-25 12:32:59", "ワードプレスサイト内に、chatgptのテキスト自動作成と自動でコピー状態の設定", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420440", 350, "🔥GPT,ChatGPT ◇サイト"),
    @("2025-10-25 12:32:59", "【時給1,600円 / 学生限定】AIでプロダクトを生成したことがある学生の方を大募集!!", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420120", 303, "🔥AI,Ai"),
    @("2025-10-25 12:32:59", "【学生発スタートアップ】留学×住まいマッチングアプリ開発仲間募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420198", 100, "◆開発 ◇アプリ"),
    @("2025-10-25 12:32:59", "wordpressレンダリングを妨げるリソースの除外", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5016989", 33, "○WordPress"),
    @("2025-10-25 12:32:59", "【クリニック向け】セキュアなPC管理の遠隔保守方法を教えてください(助言のみでのお支払い)", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420306", 30, "◇管理"),
    @("2025-10-25 12:32:59", "【Ubuntu】MySQLデータを自動CSV化しクラウド保存構築", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420180", 30, "◇MySQL"),
    @("2025-10-25 12:32:59", "【急募】時間単位で入札できるシステム構築の依頼", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5411365", 40, $null),
    @("2025-10-25 12:32:59", "【急募】monday.comとLINE WORKS連携のWebhook構築依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420427", 13, $null),
    @("2025-10-25 12:32:59", "【特急・急募】TELEC技適取得概算見積書の入手をお任せさせて頂けませんか?", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420405", 13, $null),
    @("2025-10-25 12:32:59", "運用中HPのドメイン分け", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420277", 13, $null),
    @("2025-10-25 12:32:59", "注目 【急募】YouTubeの音楽配信構築の依頼です", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420233", 13, $null),
    @("2025-10-25 12:32:59", "【急募】イベント用問い合わせLINE構築のフリーランス募集!", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5420186", 10, $null)
)

# Drop the old hyperlinks up front -- row contents are about to be fully
# rewritten below, and stale hyperlink anchors would otherwise keep pointing
# at pre-shift cells.
$ws.Range("A1:H9").Hyperlinks.Delete()

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]

    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $data[5]) | Out-Null

    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}

# スキル概要 column widened from 12 to 19 characters.
$ws.Columns.Item(8).ColumnWidth = 18.166666666666668
